$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price + 1h volume %) per the scheduled GitHub Actions refresh.
# Number-like price strings get a leading apostrophe so Excel keeps them as text
# (matching the original inline-string cells) instead of coercing them to numbers.

$ws.Range("D2").Value = '29.137.70'
$ws.Range("E2").Value = '  +1.18%  '
$ws.Range("D3").Value = '1.905.27'
$ws.Range("E3").Value = '  +1.31%  '
$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''328.16'
$ws.Range("E5").Value = '  +1.41%  '
$ws.Range("D6").Value = '''1.004'
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = '''0.4666'
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = '''0.3931'
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '''47.15'
$ws.Range("E9").Value = '  +1.39%  '
$ws.Range("D10").Value = '''0.08000'
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("D11").Value = '''1.015'
$ws.Range("E11").Value = '  +3.59%  '
$ws.Range("D12").Value = '''22.24'
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("D13").Value = '1.939.82'
$ws.Range("E13").Value = '  +2.29%  '
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("D15").Value = '''5.787'
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("D16").Value = '''0.06987'
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("D17").Value = '''89.84'
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("D19").Value = '''0.00001016'
$ws.Range("E19").Value = '  +0.60%  '
$ws.Range("D20").Value = '''17.35'
$ws.Range("E20").Value = '  +2.32%  '
$ws.Range("D21").Value = '''1.005'
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("D22").Value = '29.122.09'
$ws.Range("E22").Value = '  +1.10%  '
$ws.Range("D23").Value = '''5.376'
$ws.Range("E23").Value = '  +0.48%  '
$ws.Range("D24").Value = '''11.14'
$ws.Range("E24").Value = '  +0.43%  '
$ws.Range("D25").Value = '2.144.17'
$ws.Range("E25").Value = '  +0.99%  '
$ws.Range("D26").Value = '''2.067'
$ws.Range("E26").Value = '  -2.49%  '
$ws.Range("D27").Value = '''155.59'
$ws.Range("E27").Value = '  +1.29%  '
$ws.Range("D28").Value = '''19.79'
$ws.Range("E28").Value = '  +1.96%  '
$ws.Range("D29").Value = '''5.892'
$ws.Range("E29").Value = '  +2.30%  '
$ws.Range("D30").Value = '''2.002'
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("D31").Value = '''120.74'
$ws.Range("E31").Value = '  +0.66%  '
$ws.Range("D32").Value = '''0.09401'
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").Value = '''0.9429'
$ws.Range("E33").Value = '  +0.27%  '
$ws.Range("D34").Value = '''5.376'
$ws.Range("E34").Value = '  +1.13%  '
$ws.Range("D35").Value = '''1.359'
$ws.Range("E35").Value = '  +0.20%  '
$ws.Range("D36").Value = '''3.262'
$ws.Range("E36").Value = '  -2.49%  '
$ws.Range("D37").Value = '''0.05868'
$ws.Range("E37").Value = '  -0.81%  '
$ws.Range("D38").Value = '''1.179'
$ws.Range("E38").Value = '  +1.97%  '
$ws.Range("D39").Value = '''8.141'
$ws.Range("E39").Value = '  +3.03%  '
$ws.Range("D40").Value = '''0.02108'
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("D41").Value = '''0.5866'
$ws.Range("E41").Value = '  +2.48%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D43").Value = '''0.1821'
$ws.Range("E43").Value = '  +1.35%  '
$ws.Range("D44").Value = '''10.04'
$ws.Range("E44").Value = '  +0.39%  '
$ws.Range("D45").Value = '''2.287'
$ws.Range("E45").Value = '  +7.41%  '
$ws.Range("D46").Value = '''0.5480'
$ws.Range("E46").Value = '  +2.52%  '
$ws.Range("D47").Value = '''11.93'
$ws.Range("E47").Value = '  +0.69%  '
$ws.Range("D48").Value = '''0.07224'
$ws.Range("E48").Value = '  -1.07%  '
$ws.Range("D49").Value = '''1.883'
$ws.Range("E49").Value = '  +1.90%  '
$ws.Range("D50").Value = '''1.119'
$ws.Range("E50").Value = '  -3.86%  '
$ws.Range("D51").Value = '''113.58'
$ws.Range("E51").Value = '  -0.41%  '
